# Edit script: add columns I (I0) and J (IF) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the used range dimension implicitly happens when we write cells.

# Header row (row 1) - style matches other header cells (bold, centered, bordered) -> style index 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for rows 2..65, columns I (I0) and J (IF)
$data = @(
    "8,8",
    "9,9",
    "10,10",
    "8,8",
    "9,9",
    "9,9",
    "9,9",
    "9,9",
    "9,9",
    "8,8",
    "9,9",
    "9,9",
    "9,9",
    "8,9",
    "7,7",
    "7,7",
    "7,7",
    "7,7",
    "6,7",
    "6,7",
    "7,7",
    "7,7",
    "6,7",
    "5,5",
    "4,5",
    "8,8",
    "5,5",
    "6,7",
    "5,5",
    "7,7",
    "8,9",
    "6,6",
    "7,8",
    "8,8",
    "3,3",
    "7,7",
    "9,9",
    "4,4",
    "8,8",
    "7,8",
    "7,7",
    "7,7",
    "7,7",
    "7,7",
    "6,6",
    "6,7",
    "7,8",
    "6,6",
    "7,7",
    "7,7",
    "6,6",
    "5,5",
    "5,6",
    "6,6",
    "7,8",
    "8,8",
    "8,8",
    "7,7",
    "9,9",
    "4,4",
    "7,7",
    "7,7",
    "3,3",
    "3,3"
)

$startRow = 2
for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $parts = $data[$idx].Split(",")
    $row = $startRow + $idx
    $ws.Cells.Item($row, 9).Value = [double]$parts[0]
    $ws.Cells.Item($row, 10).Value = [double]$parts[1]
}
